$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.443.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "'1.900.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.58%  "
$ws.Range("D5").Value = "'325.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").Value = "'0.4842"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.06%  "
$ws.Range("D8").Value = "'0.4068"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.73%  "
$ws.Range("D9").Value = "'0.08091"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("D10").Value = "'1.004"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").Value = "'23.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.76%  "
$ws.Range("D12").Value = "'1.891.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").Value = "'5.976"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "'7.063"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("D15").Value = "'90.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").Value = "'1.007"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "'0.06712"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("D18").Value = "'0.00001035"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").Value = "'17.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").Value = "'1.005"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").Value = "'29.468.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "'5.564"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").Value = "'11.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.47%  "
$ws.Range("D24").Value = "'2.157"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.37%  "
$ws.Range("D25").Value = "'2.131.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").Value = "'153.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("D27").Value = "'19.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").Value = "'6.145"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.08%  "
$ws.Range("D29").Value = "'2.092"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").Value = "'118.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.70%  "
$ws.Range("D31").Value = "'1.033"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.16%  "
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("D33").Value = "'5.519"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.51%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'3.556"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.387"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.37%  "
$ws.Range("D36").Value = "'0.06090"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").Value = "'0.02256"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("D39").Value = "'0.5901"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").Value = "'7.893"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.55%  "
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").Value = "'10.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.69%  "
$ws.Range("D43").Value = "'1.294"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.69%  "
$ws.Range("D44").Value = "'0.07815"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("D45").Value = "'2.383"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.83%  "
$ws.Range("D46").Value = "'12.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.92%  "
$ws.Range("D47").Value = "'0.5534"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("D48").Value = "'1.930"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").Value = "'114.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("D50").Value = "'72.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "'1.048"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.94%  "
